# POM ChangePassword negative test and error messages
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet1: loginInfo
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("loginInfo")

# Remove the two existing hyperlinks (A2 -> buquxahu@cars2.club, A3 -> filllipa12345@yahoo.com)
$ws1.Hyperlinks.Delete()

$ws1.Range("A1").Value = "userName"
$ws1.Range("B1").Value = "passWord"
$ws1.Range("A2").Value = "filllipa12345@yahoo.com"
$ws1.Range("B2").Value = "phillipa12345"
$ws1.Range("A3").Value = "MashaRey12345@yahoo.com"
$ws1.Range("A3").Style = "Normal"
$ws1.Range("B3").Value = "masha234"

# Re-add the single hyperlink that remains (A2 -> filllipa12345@yahoo.com)
$ws1.Hyperlinks.Add($ws1.Range("A2"), "mailto:filllipa12345@yahoo.com")

$ws1.Range("B3").Select()

# ---------------------------------------------------------------------
# Sheet3: verifyRegistration -- just a selection/view change
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("verifyRegistration")
$ws3.Range("B4").Select()

# ---------------------------------------------------------------------
# Sheet4: rename the empty "Sheet2" -> "changePasswordNegative" and
# populate it with the ChangePassword negative-test data table.
# ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("Sheet2")
$ws4.Name = "changePasswordNegative"

$ws4.Cells.Item(1,1).Value = "currentPassword"
$ws4.Cells.Item(1,2).Value = "newPassword"
$ws4.Cells.Item(1,3).Value = "confirmNewPassword"
$ws4.Cells.Item(1,4).Value = "errorMessage"

$ws4.Cells.Item(2,1).Value = "phillipa12345"
$ws4.Cells.Item(2,2).Value = "milla12345"
$ws4.Cells.Item(2,3).Value = "milla123"
$ws4.Cells.Item(2,4).Value = "The passwords do not match."

$ws4.Cells.Item(3,1).Value = "phillipa12345"
$ws4.Cells.Item(3,2).Value = "mi"
$ws4.Cells.Item(3,3).Value = "mi"
$ws4.Cells.Item(3,4).Value = "This value is too short. It should have 6 characters or more."

$ws4.Cells.Item(4,1).Value = "phillipa"
$ws4.Cells.Item(4,2).Value = "kisulea"
$ws4.Cells.Item(4,3).Value = "kisulea"
$ws4.Cells.Item(4,4).Value = "Current password is wrong."

$longPassword = "1234567890123456789012345678901234567890123456789012345678901234567890123456789012345678901234567890123456789012345678901234567890123456789012345678901234567890123456789012345678901234567890123456789012345678901234567890123456789012345678901234567890123456789012345678901234567890123456789012345678901234567890123456789012345678901234567890"

$ws4.Cells.Item(5,1).Value = "phillipa12345"

$ws4.Cells.Item(5,2).NumberFormat = "@"
$ws4.Cells.Item(5,2).Value = $longPassword
$ws4.Cells.Item(5,2).Style = "Normal"

$ws4.Cells.Item(5,3).NumberFormat = "@"
$ws4.Cells.Item(5,3).Value = $longPassword
$ws4.Cells.Item(5,3).Style = "Normal"

$ws4.Cells.Item(5,4).Value = "This value is too long. It should have 255 characters or fewer."

$ws4.Columns.Item(1).AutoFit()
$ws4.Columns.Item(2).AutoFit()
$ws4.Columns.Item(3).AutoFit()
$ws4.Columns.Item(4).AutoFit()

# ---------------------------------------------------------------------
# Sheet5: brand-new "userPass" sheet, appended at the end.
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws5 = $wb.Worksheets.Add($null, $lastSheet)
$ws5.Name = "userPass"

$ws5.Range("A1").Value = "email"
$ws5.Range("B1").Value = "password"
$ws5.Range("A2").Value = "filllipa12345@yahoo.com"
$ws5.Range("B2").Value = "phillipa12345"

$ws5.Hyperlinks.Add($ws5.Range("A2"), "mailto:filllipa12345@yahoo.com")
$ws5.PageSetup.Orientation = 1

$ws5.Columns.Item(1).AutoFit()
$ws5.Columns.Item(2).AutoFit()

$ws5.Range("A2").Select()

# ---------------------------------------------------------------------
# Final selection state: changePasswordNegative (D5) should end up the
# active tab, matching activeTab="3" in the saved workbook.
# ---------------------------------------------------------------------
$ws4.Range("D5").Select()
